$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 8-17 (A=id, B=name, C, D, E=in_service).
# Two new "line" entries (line7, line8) are inserted into the name list right
# after line6, pushing extr1..extr8 down two rows. That puts line7/line8 on
# rows 8/9 (previously extr1/extr2), extr1..extr6 on rows 10-15 (previously
# extr3..extr8), and appends two brand-new rows 16/17 for extr7/extr8.

$rows = @(
  @{ Row = 8;  A = 6;  B = "line7"; C = 14; D = 11; E = $true  },
  @{ Row = 9;  A = 7;  B = "line8"; C = 16; D = 9;  E = $true  },
  @{ Row = 10; A = 8;  B = "extr1"; C = 5;  D = 12; E = $true  },
  @{ Row = 11; A = 9;  B = "extr2"; C = 5;  D = 9;  E = $true  },
  @{ Row = 12; A = 10; B = "extr3"; C = 10; D = 11; E = $true  },
  @{ Row = 13; A = 11; B = "extr4"; C = 7;  D = 8;  E = $false },
  @{ Row = 14; A = 12; B = "extr5"; C = 9;  D = 11; E = $false },
  @{ Row = 15; A = 13; B = "extr6"; C = 7;  D = 11; E = $false },
  @{ Row = 16; A = 14; B = "extr7"; C = 5;  D = 7;  E = $true  },
  @{ Row = 17; A = 15; B = "extr8"; C = 8;  D = 5;  E = $false }
)

# Rows 16 and 17 are brand new, so first clone the existing row 15's cell
# formatting into them (id column keeps the bordered/bold/centered style,
# the rest stay unstyled) before writing the new values.
$ws.Range("A15:E15").Copy()
$ws.Range("A16:E17").PasteSpecial(-4122)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
}
